# Applies crypto price/volume/coin updates captured in the commit diff.
# Each D-column price that is unambiguously parseable as a plain number is
# written with a leading apostrophe to force Excel to keep it as text (matching
# the original inlineStr cell type), then the style is reset to "Normal" so no
# stray text-format style gets attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.551.15"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.951.55"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'243.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "'0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("B7").Value = "Solana"
$ws.Range("C7").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D7").Value = "'57.70"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.68%  "
$ws.Range("D10").Value = "'0.0787"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.15%  "
$ws.Range("D11").Value = "'0.102"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("D12").Value = "2.237.62"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.822"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'21.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.46%  "
$ws.Range("D15").Value = "'13.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "'5.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").Value = "1.964.50"
$ws.Range("E17").Value = "  +0.78%  "
$ws.Range("D18").Value = "36.388.34"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "'69.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").Value = "0.0₃0844"
$ws.Range("E20").Value = "  -4.54%  "
$ws.Range("D21").Value = "'227.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'2.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("E25").Value = "  +2.44%  "
$ws.Range("D26").Value = "'9.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.138"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.27%  "
$ws.Range("D28").Value = "'160.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "'19.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.62%  "
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").Value = "'1.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "'4.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("D33").Value = "'0.0610"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.84%  "
$ws.Range("D34").Value = "'4.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.16%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("D36").Value = "'2.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.37%  "
$ws.Range("D37").Value = "'3.39"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +13.06%  "
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "'5.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -14.31%  "
$ws.Range("D40").Value = "'0.0968"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("D41").Value = "'2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").Value = "'1.15"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "'0.0209"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.10%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.359.12"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'15.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("D47").Value = "'87.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.68%  "
$ws.Range("D48").Value = "'7.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.128.33"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").Value = "'43.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.82%  "
